# Revert "Powerpoint writer: consolidate text run nodes."
#
# The original commit merged adjacent <a:r> runs that shared identical
# formatting (e.g. "A " + "slide" stayed separate, but a leading word plus
# the following space, such as "A " or "Plus ", had been written out as a
# single run). This script re-splits those runs back into a word run and a
# standalone space run, without touching any text content or formatting.
#
# Trick: writing identical text back into a TextRange.Characters(start,len)
# sub-range forces the host to re-materialize run boundaries at that
# offset (it keeps the same, empty <a:rPr/>), which is exactly the
# "un-consolidation" we need - it does not alter the visible text at all.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Slide 1, Title shape ("A slide") ---------------------------------
# "A " + "slide"  ->  "A" + " " + "slide"
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Characters(1, 1).Text = "A"

# --- Slide 1, TextBox ("Plus an image") --------------------------------
# "Plus " + "an " + "image"  ->  "Plus" + " " + "an" + " " + "image"
$textBox = $s.Shapes.Item(7)
$textBoxRange = $textBox.TextFrame.TextRange
$textBoxRange.Characters(1, 4).Text = "Plus"
$textBoxRange.Characters(6, 2).Text = "an"
